$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: copy formatting only from a source cell to a destination cell,
# then (optionally) set the destination's value.
# NOTE: this interpreter does not bind named (-Param) arguments to
# function params reliably, so use positional arguments only.
function Set-FormattedCell {
    param([string]$SrcAddr, [string]$DstAddr, $Value)
    $ws.Range($SrcAddr).Copy() | Out-Null
    $ws.Range($DstAddr).PasteSpecial(-4122) | Out-Null
    if ($null -ne $Value) {
        $ws.Range($DstAddr).Value = $Value
    }
}

# Row 3: blank separator cell, formatting only (same style as O3)
Set-FormattedCell "O3" "P3" $null

# Row 4: year header 2021
Set-FormattedCell "O4" "P4" 2021

# Row 5: total value for 2021 (style matches O8's style, per source data)
Set-FormattedCell "O8" "P5" 9038

# Rows 6-9: formatting for these four rows all uses the same cell style
# (index 17 in the stylesheet); O6 already carries that style, so reuse it
# as the format source for all four destination cells.
Set-FormattedCell "O6" "P6" $null
Set-FormattedCell "O6" "P7" 8587
Set-FormattedCell "O6" "P8" 451
Set-FormattedCell "O6" "P9" $null

# Rows 10-24: "no data" ellipsis marker, same style as the corresponding O cell
$ellipsis = [char]0x2026
for ($r = 10; $r -le 24; $r++) {
    $srcAddr = "O" + $r
    $dstAddr = "P" + $r
    Set-FormattedCell $srcAddr $dstAddr $ellipsis
}

# Row 25: bottom border row, ellipsis marker, same style as O25
Set-FormattedCell "O25" "P25" $ellipsis

# Leave the active selection on Q4, matching the post-edit state
$ws.Range("Q4").Select()
